$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 116: add the masterclass note in column F
$ws.Range("F116").Value = "Masterclass con JuanMa a las 20:00 hs."

# Row 116's date cell (C116) had lost the usual weekday fill; bring it back
# in line with the rest of the column by copying the format from C115.
$ws.Range("C115").Copy()
$ws.Range("C116").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fix the Encuentro numbering for rows 118-131 (class 17 "Front end" keeps row118
# as its own single-row entry, and a new class 18 "Sprint" spans rows 119-131)
$ws.Range("B118").Value = 1
$ws.Range("B119").Value = 1
$ws.Range("B120").Value = 2
$ws.Range("B121").Value = 3
$ws.Range("B122").Value = 4
$ws.Range("B123").Value = 5
$ws.Range("B124").Value = 6
$ws.Range("B125").Value = 7
$ws.Range("B126").Value = 8
$ws.Range("B127").Value = 9
$ws.Range("B128").Value = 10
$ws.Range("B129").Value = 11
$ws.Range("B130").Value = 12
$ws.Range("B131").Value = 13

# New class 18: "Sprint" over rows 119:131, merged in D and E columns
$ws.Range("D119:D131").Merge()
$ws.Range("D119").Value = 18

$ws.Range("E119:E131").Merge()
$ws.Range("E119").Value = "Sprint"

$ws.Range("E119:E131").Select()
